$d = $word.ActiveDocument

$replacements = @(
    @{old="643÷4=160, 3"; new="861÷9=95, 6"},
    @{old="382÷7=54, 4"; new="645÷9=71, 6"},
    @{old="774÷8=96, 6"; new="755÷6=125, 5"},
    @{old="948÷5=189, 3"; new="143÷3=47, 2"},
    @{old="951÷7=135, 6"; new="609÷9=67, 6"},
    @{old="289÷2=144, 1"; new="724÷8=90, 4"},
    @{old="788÷2=394, 0"; new="907÷3=302, 1"},
    @{old="283÷2=141, 1"; new="571÷9=63, 4"},
    @{old="466÷9=51, 7"; new="235÷6=39, 1"},
    @{old="801÷2=400, 1"; new="684÷6=114, 0"},
    @{old="446÷8=55, 6"; new="449÷4=112, 1"},
    @{old="360÷3=120, 0"; new="598÷5=119, 3"},
    @{old="457÷5=91, 2"; new="758÷7=108, 2"},
    @{old="684÷5=136, 4"; new="109÷4=27, 1"},
    @{old="817÷3=272, 1"; new="361÷6=60, 1"},
    @{old="899÷8=112, 3"; new="276÷4=69, 0"},
    @{old="666÷6=111, 0"; new="986÷6=164, 2"},
    @{old="900÷8=112, 4"; new="696÷8=87, 0"},
    @{old="261÷2=130, 1"; new="623÷2=311, 1"},
    @{old="469÷3=156, 1"; new="254÷7=36, 2"},
    @{old="249÷3=83, 0"; new="740÷9=82, 2"},
    @{old="193÷5=38, 3"; new="500÷4=125, 0"},
    @{old="483÷7=69, 0"; new="950÷8=118, 6"},
    @{old="802÷2=401, 0"; new="838÷6=139, 4"},
    @{old="806÷4=201, 2"; new="865÷7=123, 4"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
